$wb = $excel.ActiveWorkbook

# --- Sheet 1: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("L31").Value = 0
$ws1.Range("D34").Value = 2364.09
$ws1.Range("M42").Value = 1449.08
$ws1.Range("I48").Value = 157.3
$ws1.Range("E58").Value = 512.5
$ws1.Range("L58").Value = 253.44
$ws1.Range("M60").Value = "7 de 58"

# --- Sheet 2: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F31").Value = 7.13
$ws2.Range("F34").Value = 2509.16
$ws2.Range("F42").Value = 1449.08
$ws2.Range("F48").Value = 810.28
$ws2.Range("F58").Value = 739.4400000000001
$ws2.Range("F60").Value = 9436.629999999999

# --- Sheet 3: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Columns.Item(5).ColumnWidth = 22.15

$ws3.Range("D3").Value = 3752.44
$ws3.Range("E3").Value = 11072.97
$ws3.Range("F3").Value = 0.2531086829976372

$ws3.Range("D4").Value = 658.6900000000001
$ws3.Range("E4").Value = 130.6899999999999
$ws3.Range("F4").Value = 0.8344396868428388

$ws3.Range("D7").Value = 235.2
$ws3.Range("E7").Value = 651.511016287574
$ws3.Range("F7").Value = 0.2652498905277174

$ws3.Range("D11").Value = 1910.01
$ws3.Range("E11").Value = 14237.99
$ws3.Range("F11").Value = 0.1182815209313847

$ws3.Range("D12").Value = 2673.38
$ws3.Range("E12").Value = 47633.62
$ws3.Range("F12").Value = 0.05314131234221878

$ws3.Range("D14").Value = 9436.629999999999
$ws3.Range("E14").Value = 88425.25766749099
$ws3.Range("F14").Value = 0.09642803981119996
